$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-17 down to 8-18.
$ws.Rows.Item(7).Insert()

# Copy the date-format style (column D uses a custom date style) from the
# row above (row 6, which keeps its original values) into the freshly
# inserted row 7, without touching the rest of the row's formatting.
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the data for the new record (row 7).
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44453
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100114007
$ws.Cells.Item(7, 7).Value = "Jengibre"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 50
$ws.Cells.Item(7, 11).Value = 14000
$ws.Cells.Item(7, 12).Value = 15000
$ws.Cells.Item(7, 13).Value = 14600
$ws.Cells.Item(7, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 1123
$ws.Cells.Item(7, 17).Value = 13
$ws.Cells.Item(7, 18).Value = "Hortaliza"
